# ---------------------------------------------------------------------------
# Edit summary (from the target OOXML diff):
#   1) Three tables (on slides 14, 15 and 16 - each the first/only table
#      shape on its slide) switch their <a:tableStyleId> from the deck's
#      custom style {2EF2AF6F-0A52-47D6-A562-7AE053AF4BE0} to the built-in
#      style {0DEEB2D3-7035-4EE4-A8E2-CD3FC609551B}.
#   2) The presentation's theme colour scheme (ppt/theme/theme1.xml, the
#      theme actually used by the slide master / all slides) is switched
#      from the "Integral" deck theme's "Red Violet" colours to the
#      standard "Office Theme" colours.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables --------------------------------------
$newTableStyle = "{0DEEB2D3-7035-4EE4-A8E2-CD3FC609551B}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyle)
    }
}

# --- 2) Swap the theme colour scheme used by the slides -----------------
# Index -> (slot name, target "Office" RGB hex) in the order exposed by
# ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
